$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right count 4 -> 5, Wrong mark -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total): Right total 92 -> 115, Wrong mark -1 -> -1.2, Max string "91/112" -> "113.8/140"
$ws.Range("B12").Value = 115
$ws.Range("C12").Value = -1.2
$ws.Range("E12").Value = "113.8/140"
